$d = $word.ActiveDocument

$d.Content.Find.Execute("226×8=", $true, $false, $false, $false, $false, $true, 1, $false, "427×6=", 2) | Out-Null
$d.Content.Find.Execute("323×9=", $true, $false, $false, $false, $false, $true, 1, $false, "602×4=", 2) | Out-Null
$d.Content.Find.Execute("847×2=", $true, $false, $false, $false, $false, $true, 1, $false, "650×5=", 2) | Out-Null
$d.Content.Find.Execute("441×7=", $true, $false, $false, $false, $false, $true, 1, $false, "402×7=", 2) | Out-Null
$d.Content.Find.Execute("804×5=", $true, $false, $false, $false, $false, $true, 1, $false, "399×3=", 2) | Out-Null
$d.Content.Find.Execute("814×6=", $true, $false, $false, $false, $false, $true, 1, $false, "618×2=", 2) | Out-Null
$d.Content.Find.Execute("875×9=", $true, $false, $false, $false, $false, $true, 1, $false, "896×3=", 2) | Out-Null
$d.Content.Find.Execute("310×6=", $true, $false, $false, $false, $false, $true, 1, $false, "658×6=", 2) | Out-Null
$d.Content.Find.Execute("844×8=", $true, $false, $false, $false, $false, $true, 1, $false, "654×9=", 2) | Out-Null
$d.Content.Find.Execute("161×7=", $true, $false, $false, $false, $false, $true, 1, $false, "572×4=", 2) | Out-Null
$d.Content.Find.Execute("220×8=", $true, $false, $false, $false, $false, $true, 1, $false, "246×5=", 2) | Out-Null
$d.Content.Find.Execute("778×9=", $true, $false, $false, $false, $false, $true, 1, $false, "980×3=", 2) | Out-Null
$d.Content.Find.Execute("817×6=", $true, $false, $false, $false, $false, $true, 1, $false, "249×6=", 2) | Out-Null
$d.Content.Find.Execute("290×9=", $true, $false, $false, $false, $false, $true, 1, $false, "752×4=", 2) | Out-Null
$d.Content.Find.Execute("787×9=", $true, $false, $false, $false, $false, $true, 1, $false, "439×6=", 2) | Out-Null
$d.Content.Find.Execute("628×2=", $true, $false, $false, $false, $false, $true, 1, $false, "970×9=", 2) | Out-Null
$d.Content.Find.Execute("367×9=", $true, $false, $false, $false, $false, $true, 1, $false, "323×9=", 2) | Out-Null
$d.Content.Find.Execute("916×6=", $true, $false, $false, $false, $false, $true, 1, $false, "719×3=", 2) | Out-Null
$d.Content.Find.Execute("485×7=", $true, $false, $false, $false, $false, $true, 1, $false, "278×7=", 2) | Out-Null
$d.Content.Find.Execute("759×8=", $true, $false, $false, $false, $false, $true, 1, $false, "195×8=", 2) | Out-Null
$d.Content.Find.Execute("530×4=", $true, $false, $false, $false, $false, $true, 1, $false, "662×5=", 2) | Out-Null
$d.Content.Find.Execute("122×3=", $true, $false, $false, $false, $false, $true, 1, $false, "251×3=", 2) | Out-Null
$d.Content.Find.Execute("399×7=", $true, $false, $false, $false, $false, $true, 1, $false, "577×5=", 2) | Out-Null
$d.Content.Find.Execute("128×2=", $true, $false, $false, $false, $false, $true, 1, $false, "997×6=", 2) | Out-Null
$d.Content.Find.Execute("276×2=", $true, $false, $false, $false, $false, $true, 1, $false, "993×3=", 2) | Out-Null
